$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# Enter new Newspoll (row 5) region swing values
$ws.Range("B5").Value = 55
$ws.Range("C5").Value = 54
$ws.Range("D5").Value = 58
$ws.Range("E5").Value = 46
$ws.Range("F5").Value = 53
$ws.Range("G5").Value = 59

# Move the active selection to H13 as left by the author after editing
$ws.Activate()
$ws.Range("H13").Select()

$wb.Save()
